$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

# Update classification values from "critical" to "semi-critical" for the specified rows
$ws.Range("D7").Value = "semi-critical"
$ws.Range("D9").Value = "semi-critical"
$ws.Range("D14").Value = "semi-critical"
$ws.Range("D18").Value = "semi-critical"

# Update the selected cell in the sheet view
$ws.Range("H6").Select()
